$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# Step 1: locate the existing " kuch" run at the end of paragraph 2 and split
# it into a plain " " run plus a spell-checked "kuch" run, matching the
# target markup exactly. We find the word "kuch" (without the leading
# space) so the replacement range covers only the four letters.
$find = $d.Content
$find.Find.ClearFormatting()
$found = $find.Find.Execute("kuch", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$kuchRange = $d.Range($find.Start, $find.End)

$kuchXml = $pkgOpen + '<w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>kuch</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body>' + $pkgClose
$kuchRange.InsertXML($kuchXml)

# Step 2: insert the new paragraph right before the (now relocated) bookmark
# that sits at the end of paragraph 2. A collapsed range placed there causes
# InsertXML to add a brand-new paragraph immediately after the current one,
# leaving the bookmark untouched in paragraph 2.
$bm = $d.Bookmarks("_GoBack")
$ins = $d.Range($bm.Start, $bm.Start)

$paraXml = $pkgOpen + '<w:body><w:p>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Chlo</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> or </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>btao</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>kya</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>haal</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>chl</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> h&#8230;</w:t></w:r>' + `
  '</w:p></w:body>' + $pkgClose

$ins.InsertXML($paraXml)

Write-Host $d.Content.Text
